# Updated cryptos list values (price + 1h volume change) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.048.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3804"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9334"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07759"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.906.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.502"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.628"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008837"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.072.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.181"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.142.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.919"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.969"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08949"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.312"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7752"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.680"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.641"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02062"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.111"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05313"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5495"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.991"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.033"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1532"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.514"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.651"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06069"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
